# Edición del docx por Iker
#
# 1) El primer párrafo (título) pasa a tener interlineado/espaciado fijo
#    (w:spacing before/after = 0, sin auto-spacing), la marca de párrafo
#    cambia su idioma a es-ES, se añade un nuevo texto en negrita
#    "Herramientas computacionales: el arte de la programación" justo
#    después del salto de línea, y el párrafo se divide en dos: el nuevo
#    texto se queda en el párrafo del título y "Equipo 4- Sopa de
#    Ingenieros" pasa a su propio párrafo (con el mismo estilo/formato de
#    párrafo que el título).
# 2) Tres párrafos con viñetas tenían la frase final partida en dos runs
#    (texto + un run aparte solo con el punto "."). Se fusionan en un
#    único run con el texto completo terminado en punto.

$d = $word.ActiveDocument

# --- 1) Párrafo de título: spacing fijo, idiomas, nuevo subtítulo y split ---

$titlePara = $d.Paragraphs(1)
$titleXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="NormalWeb"/><w:spacing w:before="0" w:beforeAutospacing="0" w:after="0" w:afterAutospacing="0"/><w:jc w:val="center"/><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:color w:val="000000"/><w:lang w:val="es-ES"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:color w:val="000000"/><w:lang w:val="es-MX"/></w:rPr><w:t>Información solución para un problema medioambiental utilizando IoT</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:color w:val="000000"/><w:lang w:val="es-MX"/></w:rPr><w:br/></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:color w:val="000000"/><w:lang w:val="es-ES"/></w:rPr><w:t>Herramientas computacionales: el arte de la programación</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="NormalWeb"/><w:spacing w:before="0" w:beforeAutospacing="0" w:after="0" w:afterAutospacing="0"/><w:jc w:val="center"/><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:color w:val="000000"/><w:lang w:val="es-ES"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:color w:val="000000"/><w:lang w:val="es-MX"/></w:rPr><w:t>Equipo 4- Sopa de Ingenieros</w:t></w:r></w:p>
'@
$titlePara.Range.InsertXML($titleXml) | Out-Null

# --- 2) Fusionar "texto" + "." sueltos en un único run, en las tres viñetas ---

$d.Content.Find.Execute(
    "La alta propagación del COVID-19 en algunas áreas puede estar vinculada a la existencia de niveles altos de material particulado en el aire.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "La alta propagación del COVID-19 en algunas áreas puede estar vinculada a la existencia de niveles altos de material particulado en el aire.",
    2) | Out-Null

$d.Content.Find.Execute(
    "Exposición crónica a contaminantes como el dióxido de nitrógeno y el material particulado se relacionan con un aumento de la mortalidad por COVID-19.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Exposición crónica a contaminantes como el dióxido de nitrógeno y el material particulado se relacionan con un aumento de la mortalidad por COVID-19.",
    2) | Out-Null

$d.Content.Find.Execute(
    "La exposición a largo plazo a material particulado menor de 2,5 micras afecta al sistema respiratorio y cardiovascular, lo que exacerba la gravedad del COVID-19.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "La exposición a largo plazo a material particulado menor de 2,5 micras afecta al sistema respiratorio y cardiovascular, lo que exacerba la gravedad del COVID-19.",
    2) | Out-Null

Write-Host "Done. Paragraphs:" $d.Paragraphs.Count
